$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-run SGNN dialog act annotations following transcript clean up.
$ws.Range("I2").Value = "sd"
$ws.Range("J2").Value = "Statement-non-opinion"

$ws.Range("I5").Value = "aa"
$ws.Range("J5").Value = "Agree/Accept"

$ws.Range("I8").Value = "aa"
$ws.Range("J8").Value = "Agree/Accept"

$ws.Range("I28").Value = "b"
$ws.Range("J28").Value = "Acknowledge (Backchannel)"

$ws.Range("I63").Value = "sd"
$ws.Range("J63").Value = "Statement-non-opinion"

$ws.Range("I65").Value = "%"
$ws.Range("J65").Value = "Uninterpretable"
